$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.0616
$ws.Range("C3").Value = -0.1429
$ws.Range("C4").Value = -0.0374
$ws.Range("C5").Value = -0.2784
$ws.Range("C6").Value = -0.2196
